$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gruppen")

# Insert a new row before row 13 (shifts rows 13+ down by one)
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new group entry
$ws.Range("A13").Value = "Yes"
$ws.Range("B13").Value = "SecurityGroup"
$ws.Range("C13").Value = "ALYASG-ADM-LEGACYAUTHENABLED"
$ws.Range("D13").Value = "ALYASG-ADM-LEGACYAUTHENABLED"
$ws.Range("E13").Value = "Legacy Auth is not blocked by conditional Access for members in this group"
$ws.Range("G13").Value = "Private"
$ws.Range("A13:G13").Font.Bold = $true

$wb.Save()
